$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.478.21"
$ws.Range("E2").Value = "  +0.35%  "

# Row 3
$ws.Range("D3").Value = "1.838.51"

# Row 4
$ws.Range("E4").Value = "  +0.11%  "

# Row 5
$ws.Range("D5").Value = "259.88"
$ws.Range("E5").Value = "  +0.08%  "

# Row 7
$ws.Range("D7").Value = "0.5242"
$ws.Range("E7").Value = "  +0.45%  "

# Row 8
$ws.Range("D8").Value = "0.3186"
$ws.Range("E8").Value = "  -1.34%  "

# Row 9
$ws.Range("D9").Value = "0.06783"
$ws.Range("E9").Value = "  +0.29%  "

# Row 10
$ws.Range("D10").Value = "18.72"
$ws.Range("E10").Value = "  +0.65%  "

# Row 11
$ws.Range("D11").Value = "0.7834"
$ws.Range("E11").Value = "  +2.74%  "

# Row 12
$ws.Range("D12").Value = "0.07745"
$ws.Range("E12").Value = "  +0.90%  "

# Row 13
$ws.Range("D13").Value = "1.844.90"
$ws.Range("E13").Value = "  +0.32%  "

# Row 14
$ws.Range("D14").Value = "87.76"
$ws.Range("E14").Value = "  -0.80%  "

# Row 15
$ws.Range("D15").Value = "5.008"
$ws.Range("E15").Value = "  -0.12%  "

# Row 16
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  +0.10%  "

# Row 17
$ws.Range("E17").Value = "  -0.70%  "

# Row 18
$ws.Range("E18").Value = "  +0.08%  "

# Row 19
$ws.Range("E19").Value = "  +0.49%  "

# Row 20
$ws.Range("D20").Value = "26.506.73"
$ws.Range("E20").Value = "  +0.18%  "

# Row 21
$ws.Range("D21").Value = "2.069.19"
$ws.Range("E21").Value = "  -0.33%  "

# Row 22
$ws.Range("D22").Value = "4.621"
$ws.Range("E22").Value = "  +1.36%  "

# Row 23
$ws.Range("D23").Value = "5.959"
$ws.Range("E23").Value = "  +0.54%  "

# Row 24
$ws.Range("D24").Value = "9.352"
$ws.Range("E24").Value = "  -0.88%  "

# Row 25
$ws.Range("D25").Value = "141.69"
$ws.Range("E25").Value = "  -2.17%  "

# Row 26
$ws.Range("D26").Value = "2.178"
$ws.Range("E26").Value = "  -2.05%  "

# Row 27
$ws.Range("D27").Value = "1.678"
$ws.Range("E27").Value = "  +1.30%  "

# Row 28
$ws.Range("E28").Value = "  -0.28%  "

# Row 29
$ws.Range("D29").Value = "111.37"
$ws.Range("E29").Value = "  +0.24%  "

# Row 30
$ws.Range("D30").Value = "4.156"
$ws.Range("E30").Value = "  -0.34%  "

# Row 31
$ws.Range("D31").Value = "0.08692"
$ws.Range("E31").Value = "  -0.58%  "

# Row 32
$ws.Range("D32").Value = "4.068"
$ws.Range("E32").Value = "  -1.63%  "

# Row 33
$ws.Range("D33").Value = "0.04877"
$ws.Range("E33").Value = "  +1.22%  "

# Row 34
$ws.Range("D34").Value = "0.7271"
$ws.Range("E34").Value = "  +3.87%  "

# Row 35
$ws.Range("D35").Value = "1.133"
$ws.Range("E35").Value = "  +1.22%  "

# Row 36
$ws.Range("D36").Value = "2.875"
$ws.Range("E36").Value = "  +1.02%  "

# Row 37
$ws.Range("D37").Value = "3.093"
$ws.Range("E37").Value = "  +1.20%  "

# Row 38
$ws.Range("D38").Value = "2.235"
$ws.Range("E38").Value = "  +2.15%  "

# Row 39
$ws.Range("D39").Value = "0.01752"
$ws.Range("E39").Value = "  -0.45%  "

# Row 40
$ws.Range("D40").Value = "0.4764"
$ws.Range("E40").Value = "  -1.28%  "

# Row 41
$ws.Range("D41").Value = "0.8920"
$ws.Range("E41").Value = "  +0.47%  "

# Row 42
$ws.Range("E42").Value = "  -1.65%  "

# Row 43
$ws.Range("D43").Value = "5.935"
$ws.Range("E43").Value = "  -2.50%  "

# Row 45
$ws.Range("D45").Value = "7.655"
$ws.Range("E45").Value = "  +0.52%  "

# Row 46
$ws.Range("D46").Value = "0.4160"
$ws.Range("E46").Value = "  +1.28%  "

# Row 47
$ws.Range("D47").Value = "8.995"
$ws.Range("E47").Value = "  +0.03%  "

# Row 48
$ws.Range("D48").Value = "0.05849"
$ws.Range("E48").Value = "  -0.05%  "

# Row 49
$ws.Range("D49").Value = "0.1230"
$ws.Range("E49").Value = "  +1.10%  "

# Row 50
$ws.Range("D50").Value = "34.82"
$ws.Range("E50").Value = "  +0.42%  "

# Row 51
$ws.Range("D51").Value = "0.8908"
$ws.Range("E51").Value = "  +1.15%  "
